$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Average Rating"
$ws.Range("F2").Formula = "=AVERAGE(B2:E2)"
